$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as literal text in the workbook
# (e.g. "67.488.11", "0.999"). Plain-decimal-looking values would be
# auto-converted to numbers by Excel on assignment, silently dropping
# trailing zeros ("1.00" -> 1) or changing precision, so force those
# specific cells to Text format first to preserve the exact string.
# Values that already contain multiple "." (thousands-grouped, e.g.
# "67.600.16") are never valid numbers, so Excel keeps them as text
# without any extra help.
$ws.Range("D2").Value = '67.600.16'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '3.763.24'
$ws.Range("E3").Value = '  -1.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.45'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.32'
$ws.Range("E6").Value = '  +2.28%  '
$ws.Range("D7").Value = '3.763.12'
$ws.Range("E7").Value = '  -1.58%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("E12").Value = '  -0.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000272'
$ws.Range("E13").Value = '  +5.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.73'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '4.394.20'
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").Value = '3.761.86'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.88'
$ws.Range("E17").Value = '  +4.21%  '
$ws.Range("D18").Value = '67.594.43'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.23'
$ws.Range("E19").Value = '  -1.51%  '
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.55'
$ws.Range("E21").Value = '  -3.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '470.11'
$ws.Range("E22").Value = '  +1.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000149'
$ws.Range("E24").Value = '  -6.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.87'
$ws.Range("E25").Value = '  +1.31%  '
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.20'
$ws.Range("E27").Value = '  +1.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.32'
$ws.Range("E28").Value = '  +3.57%  '
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("E30").Value = '  -1.83%  '
$ws.Range("D31").Value = '3.907.40'
$ws.Range("E31").Value = '  -1.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.69'
$ws.Range("E32").Value = '  +1.75%  '
$ws.Range("E33").Value = '  -1.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.45'
$ws.Range("E34").Value = '  -1.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.15'
$ws.Range("E35").Value = '  -3.80%  '
$ws.Range("D36").Value = '3.728.42'
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.85'
$ws.Range("E37").Value = '  +8.12%  '
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.90'
$ws.Range("E39").Value = '  +0.54%  '
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -1.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.71'
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.90'
$ws.Range("E47").Value = '  -2.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '398.55'
$ws.Range("E48").Value = '  -4.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000270'
$ws.Range("E49").Value = '  -8.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.45'
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0354'
$ws.Range("E51").Value = '  -0.56%  '
